$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the name "홍길동" from B2 to C3
$ws.Range("B2").Value = $null
$ws.Range("C3").Value = "홍길동"
